# Update the practice table of two-digit ÷ one-digit division problems
# to the newly generated answer set, cell by cell, preserving formatting.
$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$cell = $table.Cell(1, 1)
$cell.Range.Text = "83÷8=10, 3"
$cell = $table.Cell(1, 2)
$cell.Range.Text = "85÷2=42, 1"
$cell = $table.Cell(1, 3)
$cell.Range.Text = "83÷6=13, 5"
$cell = $table.Cell(1, 4)
$cell.Range.Text = "15÷3=5, 0"
$cell = $table.Cell(1, 5)
$cell.Range.Text = "29÷2=14, 1"

$cell = $table.Cell(5, 1)
$cell.Range.Text = "10÷3=3, 1"
$cell = $table.Cell(5, 2)
$cell.Range.Text = "74÷9=8, 2"
$cell = $table.Cell(5, 3)
$cell.Range.Text = "14÷4=3, 2"
$cell = $table.Cell(5, 4)
$cell.Range.Text = "75÷6=12, 3"
$cell = $table.Cell(5, 5)
$cell.Range.Text = "14÷9=1, 5"

$cell = $table.Cell(9, 1)
$cell.Range.Text = "33÷8=4, 1"
$cell = $table.Cell(9, 2)
$cell.Range.Text = "69÷8=8, 5"
$cell = $table.Cell(9, 3)
$cell.Range.Text = "87÷6=14, 3"
$cell = $table.Cell(9, 4)
$cell.Range.Text = "51÷4=12, 3"
$cell = $table.Cell(9, 5)
$cell.Range.Text = "80÷7=11, 3"

$cell = $table.Cell(13, 1)
$cell.Range.Text = "36÷4=9, 0"
$cell = $table.Cell(13, 2)
$cell.Range.Text = "60÷2=30, 0"
$cell = $table.Cell(13, 3)
$cell.Range.Text = "54÷5=10, 4"
$cell = $table.Cell(13, 4)
$cell.Range.Text = "55÷5=11, 0"
$cell = $table.Cell(13, 5)
$cell.Range.Text = "65÷9=7, 2"

$cell = $table.Cell(17, 1)
$cell.Range.Text = "25÷6=4, 1"
$cell = $table.Cell(17, 2)
$cell.Range.Text = "43÷3=14, 1"
$cell = $table.Cell(17, 3)
$cell.Range.Text = "15÷4=3, 3"
$cell = $table.Cell(17, 4)
$cell.Range.Text = "77÷3=25, 2"
$cell = $table.Cell(17, 5)
$cell.Range.Text = "96÷9=10, 6"
